$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update masthead text (volume/issue number and report date range) ---
# A8 shared string: "Volume 31   Number  21" -> "...22" (last run, chars 21-22)
$ws.Range("A8").Characters(21,2).Text = "22"

# C9 shared string: date range updated (first date same length, second date shorter)
$ws.Range("C9").Characters(27,9).Text = "5/27/2024"
$ws.Range("C9").Characters(47,9).Text = "6/2/2024"

# --- Style-transplant donors for cells that switch between text-placeholder and numeric styles ---
$styleDonor14 = $ws.Range("C23")   # style 14 (text placeholder style)
$styleDonor15 = $ws.Range("I14")   # style 15 (plain integer numeric style)
$styleDonor16 = $ws.Range("K14")   # style 16 (percent-like numeric style)
$starDonor    = $ws.Range("E23")   # style 14, text "***.*"

# --- Crime statistics table (rows 15-28, 31, 33) ---

# Row 15
$styleDonor14.Copy($ws.Range("C15"))
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 27.272727272727
$ws.Range("M15").Value = 600
$ws.Range("N15").Value = 7.692307692307

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = -18.75
$ws.Range("I16").Value = 151
$ws.Range("J16").Value = 208
$ws.Range("K16").Value = -27.403846153846
$ws.Range("L16").Value = -36.016949152542
$ws.Range("M16").Value = 147.540983606557
$ws.Range("N16").Value = -85.673624288425

# Row 17
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 40
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 81.818181818181
$ws.Range("I17").Value = 230
$ws.Range("J17").Value = 207
$ws.Range("K17").Value = 11.111111111111
$ws.Range("L17").Value = 50.32679738562
$ws.Range("M17").Value = 177.10843373494
$ws.Range("N17").Value = -16.96750902527

# Row 18
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 14.285714285714
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -7.142857142857
$ws.Range("I18").Value = 164
$ws.Range("J18").Value = 193
$ws.Range("K18").Value = -15.025906735751
$ws.Range("L18").Value = -41.637010676156
$ws.Range("M18").Value = 18.840579710144
$ws.Range("N18").Value = -85.460992907801

# Row 19
$ws.Range("C19").Value = 34
$ws.Range("D19").Value = 36
$ws.Range("E19").Value = -5.555555555555
$ws.Range("F19").Value = 155
$ws.Range("G19").Value = 193
$ws.Range("H19").Value = -19.689119170984
$ws.Range("I19").Value = 873
$ws.Range("J19").Value = 1028
$ws.Range("K19").Value = -15.077821011673
$ws.Range("L19").Value = 3.558718861209
$ws.Range("M19").Value = -8.394543546694
$ws.Range("N19").Value = -77.769289533995

# Row 20
$styleDonor15.Copy($ws.Range("C20"))
$ws.Range("C20").Value = 2
$styleDonor15.Copy($ws.Range("D20"))
$ws.Range("D20").Value = 1
$styleDonor16.Copy($ws.Range("E20"))
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 23
$ws.Range("J20").Value = 26
$ws.Range("K20").Value = -11.538461538461
$ws.Range("L20").Value = -23.333333333333
$ws.Range("M20").Value = 155.555555555556
$ws.Range("N20").Value = -86.309523809523

# Row 21
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 56
$ws.Range("E21").Value = 3.571428571428
$ws.Range("F21").Value = 255
$ws.Range("G21").Value = 284
$ws.Range("H21").Value = -10.211267605633
$ws.Range("I21").Value = 1457
$ws.Range("J21").Value = 1670
$ws.Range("K21").Value = -12.754491017964
$ws.Range("L21").Value = -6.542655548428
$ws.Range("M21").Value = 16.934189406099
$ws.Range("N21").Value = -77.833561539631

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = -52.380952380952
$ws.Range("I22").Value = 73
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = -27
$ws.Range("L22").Value = -9.876543209876
$ws.Range("M22").Value = 35.185185185185

# Row 24
$ws.Range("C24").Value = 91
$ws.Range("D24").Value = 82
$ws.Range("E24").Value = 10.975609756097
$ws.Range("F24").Value = 359
$ws.Range("G24").Value = 329
$ws.Range("H24").Value = 9.118541033434
$ws.Range("I24").Value = 1792
$ws.Range("J24").Value = 1589
$ws.Range("K24").Value = 12.775330396475
$ws.Range("L24").Value = 45.336577453365
$ws.Range("M24").Value = -7.102125453602

# Row 25
$ws.Range("C25").Value = 82
$ws.Range("D25").Value = 74
$ws.Range("E25").Value = 10.81081081081
$ws.Range("F25").Value = 318
$ws.Range("G25").Value = 317
$ws.Range("H25").Value = 0.315457413249
$ws.Range("I25").Value = 1580
$ws.Range("J25").Value = 1495
$ws.Range("K25").Value = 5.685618729096
$ws.Range("L25").Value = 32.773109243697

# Row 26
$ws.Range("C26").Value = 31
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 138.461538461538
$ws.Range("F26").Value = 89
$ws.Range("G26").Value = 68
$ws.Range("H26").Value = 30.882352941176
$ws.Range("I26").Value = 431
$ws.Range("J26").Value = 408
$ws.Range("K26").Value = 5.63725490196
$ws.Range("L26").Value = 21.408450704225
$ws.Range("M26").Value = 69.019607843137

# Row 27
$styleDonor14.Copy($ws.Range("C27"))
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -71.428571428571
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = 38.461538461538
$ws.Range("L27").Value = 38.461538461538

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = -6.25
$ws.Range("I28").Value = 86
$ws.Range("J28").Value = 88
$ws.Range("K28").Value = -2.272727272727
$ws.Range("L28").Value = 3.614457831325

# Row 31
$styleDonor14.Copy($ws.Range("D31"))
$starDonor.Copy($ws.Range("E31"))
$ws.Range("I31").Value = 8
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = -33.333333333333

# Row 33
$styleDonor15.Copy($ws.Range("C33"))
$ws.Range("C33").Value = 1
$styleDonor15.Copy($ws.Range("F33"))
$ws.Range("F33").Value = 1
$styleDonor15.Copy($ws.Range("I33"))
$ws.Range("I33").Value = 1
